# Append the new price-data row (row 80) to the end of the sheet,
# matching the source data's convention of storing the date as a
# text string (e.g. "2024-10-21 00:00:00") in column A and plain
# numbers in columns B-E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A80").Value = "2024-10-21 00:00:00"
$ws.Range("B80").Value = 73500
$ws.Range("C80").Value = 10331.89
$ws.Range("D80").Value = 9143.26
$ws.Range("E80").Value = 7.1132
